$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 31   Number  43"
$ws.Range("C9").Value = "Report Covering the Week  10/21/2024  Through  10/27/2024"

# --- Type-changing cells: copy format+value from a STABLE source cell with the exact target style/value, then overwrite value if needed ---
$ws.Range("C14").Copy($ws.Range("D17"))   # D17 -> text "0" (style 13)
$ws.Range("E14").Copy($ws.Range("E17"))   # E17 -> text "***.*" (style 13)
$ws.Range("G15").Copy($ws.Range("C20"))   # C20 -> numeric (style 15); value set below
$ws.Range("D14").Copy($ws.Range("C23"))   # C23 -> text "0" (style 13)
$ws.Range("F14").Copy($ws.Range("D23"))   # D23 -> text "0" (style 13)
$ws.Range("H14").Copy($ws.Range("E23"))   # E23 -> text "***.*" (style 13)
$ws.Range("C28").Copy($ws.Range("D28"))   # D28 -> numeric value 1 (style 15)
$ws.Range("K23").Copy($ws.Range("E28"))   # E28 -> numeric value 0 (style 14)

# --- Pure numeric value updates ---
$ws.Range("C20").Value = 3
$ws.Range("N15").Value = -12.5
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = -71.428571428571
$ws.Range("G16").Value = 14
$ws.Range("H16").Value = -78.571428571428
$ws.Range("I16").Value = 80
$ws.Range("J16").Value = 106
$ws.Range("K16").Value = -24.528301886792
$ws.Range("L16").Value = -32.773109243697
$ws.Range("M16").Value = -29.203539823008
$ws.Range("N16").Value = -84.282907662082
$ws.Range("C17").Value = 1
$ws.Range("F17").Value = 8
$ws.Range("G17").Value = 5
$ws.Range("H17").Value = 60
$ws.Range("I17").Value = 110
$ws.Range("K17").Value = 15.78947368421
$ws.Range("L17").Value = -16.666666666666
$ws.Range("M17").Value = 71.875
$ws.Range("N17").Value = -50.67264573991
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -75
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = 18.181818181818
$ws.Range("I18").Value = 127
$ws.Range("J18").Value = 146
$ws.Range("K18").Value = -13.013698630137
$ws.Range("L18").Value = -41.203703703703
$ws.Range("M18").Value = -38.349514563106
$ws.Range("N18").Value = -86.715481171548
$ws.Range("C19").Value = 20
$ws.Range("D19").Value = 19
$ws.Range("E19").Value = 5.263157894736
$ws.Range("F19").Value = 62
$ws.Range("G19").Value = 64
$ws.Range("H19").Value = -3.125
$ws.Range("I19").Value = 542
$ws.Range("J19").Value = 568
$ws.Range("K19").Value = -4.577464788732
$ws.Range("L19").Value = 10.386965376782
$ws.Range("M19").Value = 115.93625498008
$ws.Range("N19").Value = 77.704918032786
$ws.Range("E20").Value = 50
$ws.Range("F20").Value = 11
$ws.Range("H20").Value = 37.5
$ws.Range("I20").Value = 97
$ws.Range("J20").Value = 131
$ws.Range("K20").Value = -25.954198473282
$ws.Range("L20").Value = -36.184210526315
$ws.Range("M20").Value = -21.774193548387
$ws.Range("N20").Value = -87.032085561497
$ws.Range("C21").Value = 27
$ws.Range("D21").Value = 32
$ws.Range("E21").Value = -15.625
$ws.Range("F21").Value = 97
$ws.Range("G21").Value = 103
$ws.Range("H21").Value = -5.825242718446
$ws.Range("I21").Value = 963
$ws.Range("J21").Value = 1052
$ws.Range("K21").Value = -8.460076045627
$ws.Range("L21").Value = -14.171122994652
$ws.Range("M21").Value = 26.877470355731
$ws.Range("N21").Value = -65.058055152394
$ws.Range("I22").Value = 12
$ws.Range("K22").Value = 33.333333333333
$ws.Range("L22").Value = 9.090909090909
$ws.Range("M22").Value = -7.692307692307
$ws.Range("L23").Value = -14.285714285714
$ws.Range("M23").Value = 9.090909090909
$ws.Range("C24").Value = 34
$ws.Range("D24").Value = 19
$ws.Range("E24").Value = 78.947368421052
$ws.Range("F24").Value = 89
$ws.Range("G24").Value = 58
$ws.Range("H24").Value = 53.448275862069
$ws.Range("I24").Value = 826
$ws.Range("J24").Value = 748
$ws.Range("K24").Value = 10.427807486631
$ws.Range("L24").Value = -0.601684717208
$ws.Range("M24").Value = 59.459459459459
$ws.Range("C25").Value = 23
$ws.Range("D25").Value = 12
$ws.Range("E25").Value = 91.666666666666
$ws.Range("F25").Value = 56
$ws.Range("G25").Value = 33
$ws.Range("H25").Value = 69.696969696969
$ws.Range("I25").Value = 514
$ws.Range("J25").Value = 403
$ws.Range("K25").Value = 27.543424317617
$ws.Range("L25").Value = 9.361702127659
$ws.Range("C26").Value = 9
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = 28.571428571428
$ws.Range("F26").Value = 26
$ws.Range("G26").Value = 21
$ws.Range("H26").Value = 23.809523809523
$ws.Range("I26").Value = 227
$ws.Range("J26").Value = 204
$ws.Range("K26").Value = 11.274509803921
$ws.Range("L26").Value = 1.339285714285
$ws.Range("M26").Value = 28.248587570621
$ws.Range("F28").Value = 3
$ws.Range("H28").Value = -25
$ws.Range("I28").Value = 33
$ws.Range("J28").Value = 37
$ws.Range("K28").Value = -10.81081081081
$ws.Range("L28").Value = 32
